$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The five observation rows (17-21) on the "Artfynd" sheet got their record
# identifiers (column A) reshuffled among each other, their coordinates
# (Q = Ost / R = Nord) were rounded off to whole meters, the Starttid/Sluttid
# (Z/AB) time stamps were dropped everywhere, and the Enhet/Ålder-Stadium/
# Kön/Metod/Bestämningsmetod columns (J/K/L/N/AF) - which only applied to
# the two "blomning" (flowering) observations - moved along with their
# corresponding record identifiers.
# ---------------------------------------------------------------------------

# Helper: make a truly blank *text* cell (mirrors the original empty
# inlineStr cells such as J17/L17/N17/AF17 in the source workbook).
function Set-BlankTextCell($range) {
    $range.Value = "'"
}

# --- Row 17: Id 111821924 -> 111821926 (coords of the old row 20 record) ---
$ws.Range("A17").Value = 111821926
$ws.Range("J17").ClearContents()
$ws.Range("K17").ClearContents()
$ws.Range("L17").ClearContents()
$ws.Range("N17").ClearContents()
$ws.Range("Q17").Value = 550846
$ws.Range("R17").Value = 6681625
$ws.Range("Z17").ClearContents()
$ws.Range("AB17").ClearContents()
$ws.Range("AF17").ClearContents()

# --- Row 18: Id 111821923 -> 111821927 (coords of the old row 19 record) ---
$ws.Range("A18").Value = 111821927
$ws.Range("J18").ClearContents()
$ws.Range("K18").ClearContents()
$ws.Range("L18").ClearContents()
$ws.Range("N18").ClearContents()
$ws.Range("Q18").Value = 550820
$ws.Range("R18").Value = 6681733
$ws.Range("Z18").ClearContents()
$ws.Range("AB18").ClearContents()
$ws.Range("AF18").ClearContents()

# --- Row 19: Id 111821927 -> 111821924 (coords of the old row 17 record) ---
$ws.Range("A19").Value = 111821924
Set-BlankTextCell $ws.Range("J19")
$ws.Range("K19").Value = "blomning"
Set-BlankTextCell $ws.Range("L19")
Set-BlankTextCell $ws.Range("N19")
$ws.Range("Q19").Value = 550675
$ws.Range("R19").Value = 6681937
$ws.Range("Z19").ClearContents()
$ws.Range("AB19").ClearContents()
Set-BlankTextCell $ws.Range("AF19")

# --- Row 20: Id 111821926 -> 111821928 (coords of the old row 21 record) ---
$ws.Range("A20").Value = 111821928
$ws.Range("Q20").Value = 550826
$ws.Range("R20").Value = 6681726
$ws.Range("Z20").ClearContents()
$ws.Range("AB20").ClearContents()

# --- Row 21: Id 111821928 -> 111821923 (coords of the old row 18 record) ---
$ws.Range("A21").Value = 111821923
Set-BlankTextCell $ws.Range("J21")
$ws.Range("K21").Value = "blomning"
Set-BlankTextCell $ws.Range("L21")
Set-BlankTextCell $ws.Range("N21")
$ws.Range("Q21").Value = 550701
$ws.Range("R21").Value = 6681909
$ws.Range("Z21").ClearContents()
$ws.Range("AB21").ClearContents()
Set-BlankTextCell $ws.Range("AF21")
